# ---------------------------------------------------------------------------
# This script reproduces the structural edits described by the commit
# "Structure changes, improve documentation." It:
#
#  1. Re-labels / re-orders the metadata rows on the HEADER sheet.
#  2. Inserts a new leading "-" entry into the lookup lists that live on the
#     various hidden "<NAME>_" helper sheets (shifting the existing entries
#     down by one row), and
#  3. Updates the data-validation Formula1 ranges on the corresponding
#     visible sheets so they keep pointing at the (now one-row-larger)
#     lookup ranges.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. HEADER sheet - re-label rows 3-6 (values in column B move along with
#    the row they are already on; rows 5 & 6 lose their old values because
#    the labels that used to own them - CATEGORY/SUB_CATEGORY - are gone).
# ---------------------------------------------------------------------------
$header = $wb.Worksheets.Item("HEADER")
$header.Range("A3").Value = "DOMAIN"
$header.Range("A4").Value = "CATEGORY"
$header.Range("A5").Value = "SOURCE_ORG"
$header.Range("B5").Value = ""
$header.Range("A6").Value = "SOURCE_PERSON"
$header.Range("B6").Value = ""

# ---------------------------------------------------------------------------
# 2. Helper function: writes an array of values into column $col of sheet
#    $sheetName, starting at row 1, prefixed with a leading "-" row (i.e.
#    the existing list shifts down by one row).
# ---------------------------------------------------------------------------
function Set-ShiftedList($sheetName, $col, $values) {
    $ws = $wb.Worksheets.Item($sheetName)
    $shifted = @("-") + $values
    for ($i = 0; $i -lt $shifted.Length; $i++) {
        $ws.Cells.Item($i + 1, $col).Value = $shifted[$i]
    }
}

# The common 34-entry lookup list shared by most of the hidden helper sheets.
$list1 = @("A","B","D","H","M","N","Q","S","W","A2","A3","A4","A5","A10","A20", `
           "A30","A_3","M2","M_2","M_3","W2","W3","W4","W_2","W_3","D_2","H2", `
           "H3","I","OA","OM","_O","_U","_Z")

# The 4-entry "hedge" list used only by EXT_DERIVATIVE_ / OTC_DERIVATIVE_.
$list2 = @("MICRO_HEDGE","MACRO_HEDGE","PORT_MGMNT","SPECULATIVE")

# The 13-entry list used only by REM_FIXED_ (column A there).
$list3 = @("N1131","N1132","N1139","N114","N115","N1171","N1172","N1173", `
           "N1174","N1179","N131","N132","N133")

# Sheets whose column A holds $list1 (and nothing else changes there).
$simpleSheets = @("SHARE_","ASSET_DEBT_","RESIDENTIAL_RE_","COMMERCIAL_RE_", `
                   "LIAB_DEBT_","HOLDER_")
foreach ($name in $simpleSheets) {
    Set-ShiftedList $name 1 $list1
}

# Sheets whose column A holds $list1 AND column B holds $list2.
$hedgeSheets = @("EXT_DERIVATIVE_","OTC_DERIVATIVE_")
foreach ($name in $hedgeSheets) {
    Set-ShiftedList $name 1 $list1
    Set-ShiftedList $name 2 $list2
}

# REM_FIXED_ : column A holds $list3, column C holds $list1.
Set-ShiftedList "REM_FIXED_" 1 $list3
Set-ShiftedList "REM_FIXED_" 3 $list1

# ---------------------------------------------------------------------------
# 3. Fix up the data validation formulas on the visible sheets so the
#    referenced ranges keep up with the newly-grown lookup lists.
# ---------------------------------------------------------------------------
function Set-ValidationFormula($sheetName, $sqref, $formula) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($sqref).Validation.Formula1 = $formula
}

Set-ValidationFormula "SHARE"          "D4:D20" "'SHARE_'!`$A`$1:`$A`$35"

Set-ValidationFormula "EXT_DERIVATIVE" "D4:D20" "'EXT_DERIVATIVE_'!`$A`$1:`$A`$35"
Set-ValidationFormula "EXT_DERIVATIVE" "J4:J20" "'EXT_DERIVATIVE_'!`$B`$1:`$B`$5"

Set-ValidationFormula "OTC_DERIVATIVE" "C4:C20" "'OTC_DERIVATIVE_'!`$A`$1:`$A`$35"
Set-ValidationFormula "OTC_DERIVATIVE" "G4:G20" "'OTC_DERIVATIVE_'!`$B`$1:`$B`$5"

Set-ValidationFormula "RESIDENTIAL_RE" "C4:C20" "'RESIDENTIAL_RE_'!`$A`$1:`$A`$35"

Set-ValidationFormula "COMMERCIAL_RE"  "C4:C20" "'COMMERCIAL_RE_'!`$A`$1:`$A`$35"

Set-ValidationFormula "REM_FIXED"      "A4:A20" "'REM_FIXED_'!`$A`$1:`$A`$14"
Set-ValidationFormula "REM_FIXED"      "E4:E20" "'REM_FIXED_'!`$C`$1:`$C`$35"

Set-ValidationFormula "LIAB_DEBT"      "D4:D20" "'LIAB_DEBT_'!`$A`$1:`$A`$35"

Set-ValidationFormula "HOLDER"         "D4:D20" "'HOLDER_'!`$A`$1:`$A`$35"

Set-ValidationFormula "ASSET_DEBT"     "D4:D20" "'ASSET_DEBT_'!`$A`$1:`$A`$35"
